$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.277.35"
$ws.Range("E2").Value = "  +5.47%  "

$ws.Range("D3").Value = "2.769.72"
$ws.Range("E3").Value = "  +4.54%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.77%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.85%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.14%  "

$ws.Range("E8").Value = "  +1.77%  "

$ws.Range("D9").Value = "2.769.03"
$ws.Range("E9").Value = "  +3.81%  "

$ws.Range("E10").Value = "  +2.40%  "

$ws.Range("E11").Value = "  +5.45%  "

$ws.Range("E12").Value = "  +2.92%  "

$ws.Range("E13").Value = "  +3.17%  "

$ws.Range("D14").Value = "3.261.07"
$ws.Range("E14").Value = "  +4.39%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.84"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.88%  "

$ws.Range("D16").Value = "64.187.76"
$ws.Range("E16").Value = "  +5.29%  "

$ws.Range("E17").Value = "  +6.98%  "

$ws.Range("D18").Value = "2.769.17"
$ws.Range("E18").Value = "  +4.12%  "

$ws.Range("E19").Value = "  +3.26%  "

$ws.Range("E20").Value = "  +3.99%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "362.41"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.63%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.89%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.32%  "

$ws.Range("E24").Value = "  +0.57%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.69"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.99%  "

$ws.Range("E26").Value = "  +5.59%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.57"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.11%  "

$ws.Range("E28").Value = "  +0.33%  "

$ws.Range("D29").Value = "0.0₃0904"
$ws.Range("E29").Value = "  +11.47%  "

$ws.Range("E30").Value = "  +1.66%  "

$ws.Range("E31").Value = "  +4.45%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +21.44%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "173.36"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.83%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.06%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "20.56"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.30%  "

$ws.Range("E36").Value = "  +7.75%  "

$ws.Range("E37").Value = "  +8.44%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.84"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +10.41%  "

$ws.Range("E39").Value = "  +12.40%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "345.03"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.15%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.25"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.39%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.34"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.99%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.84"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +12.34%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.02"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.11%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "22.04"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.43%  "

$ws.Range("E46").Value = "  +6.20%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.650"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.74%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "138.23"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.18%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0256"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.38%  "

$ws.Range("E50").Value = "  +1.86%  "

$ws.Range("E51").Value = "  +0.25%  "
